$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 22 de Julio de 2020 a las 14:16"

# Row 25
$ws.Cells.Item(25,2).Value = 107871
$ws.Cells.Item(25,3).Value = 441
$ws.Cells.Item(25,4).Value = 104641
$ws.Cells.Item(25,5).Value = 3067
$ws.Cells.Item(25,6).Value = 0
$ws.Cells.Item(25,7).Value = 3
$ws.Cells.Item(25,8).Value = 163

# Row 35
$ws.Cells.Item(35,2).Value = 66521
$ws.Cells.Item(35,3).Value = 173
$ws.Cells.Item(35,4).Value = 59061
$ws.Cells.Item(35,5).Value = 6947
$ws.Cells.Item(35,6).Value = 0
$ws.Cells.Item(35,7).Value = 6
$ws.Cells.Item(35,8).Value = 513

# Row 55
$ws.Cells.Item(55,2).Value = 33883
$ws.Cells.Item(55,3).Value = 141
$ws.Cells.Item(55,4).Value = 30500
$ws.Cells.Item(55,5).Value = 1411
$ws.Cells.Item(55,6).Value = 0
$ws.Cells.Item(55,7).Value = 0
$ws.Cells.Item(55,8).Value = 1972

# Row 73
$ws.Cells.Item(73,2).Value = 13350
$ws.Cells.Item(73,3).Value = 48
$ws.Cells.Item(73,4).Value = 12274
$ws.Cells.Item(73,5).Value = 465
$ws.Cells.Item(73,6).Value = 0
$ws.Cells.Item(73,7).Value = 0
$ws.Cells.Item(73,8).Value = 611

# Row 83 (country -> Bosnia y Herzegovina)
$ws.Cells.Item(83,1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(83,2).Value = 9115
$ws.Cells.Item(83,3).Value = 328
$ws.Cells.Item(83,4).Value = 4159
$ws.Cells.Item(83,5).Value = 4688
$ws.Cells.Item(83,6).Value = 0
$ws.Cells.Item(83,7).Value = 4
$ws.Cells.Item(83,8).Value = 268

# Row 84 (country -> Noruega)
$ws.Cells.Item(84,1).Value = "Noruega"
$ws.Cells.Item(84,2).Value = 9053
$ws.Cells.Item(84,3).Value = 0
$ws.Cells.Item(84,4).Value = 8138
$ws.Cells.Item(84,5).Value = 660
$ws.Cells.Item(84,6).Value = 0
$ws.Cells.Item(84,7).Value = 0
$ws.Cells.Item(84,8).Value = 255

# Row 85 (country -> Senegal)
$ws.Cells.Item(85,1).Value = "Senegal"
$ws.Cells.Item(85,2).Value = 8985
$ws.Cells.Item(85,3).Value = 0
$ws.Cells.Item(85,4).Value = 6044
$ws.Cells.Item(85,5).Value = 2767
$ws.Cells.Item(85,6).Value = 0
$ws.Cells.Item(85,7).Value = 0
$ws.Cells.Item(85,8).Value = 174

# Row 86 (country -> Malasia)
$ws.Cells.Item(86,1).Value = "Malasia"
$ws.Cells.Item(86,2).Value = 8831
$ws.Cells.Item(86,3).Value = 16
$ws.Cells.Item(86,4).Value = 8566
$ws.Cells.Item(86,5).Value = 142
$ws.Cells.Item(86,6).Value = 0
$ws.Cells.Item(86,7).Value = 0
$ws.Cells.Item(86,8).Value = 123

# Row 88
$ws.Cells.Item(88,2).Value = 8162
$ws.Cells.Item(88,3).Value = 614
$ws.Cells.Item(88,4).Value = 4662
$ws.Cells.Item(88,5).Value = 3431
$ws.Cells.Item(88,6).Value = 0
$ws.Cells.Item(88,7).Value = 4
$ws.Cells.Item(88,8).Value = 69

# Row 99
$ws.Cells.Item(99,2).Value = 4530
$ws.Cells.Item(99,3).Value = 108
$ws.Cells.Item(99,4).Value = 3278
$ws.Cells.Item(99,5).Value = 1127
$ws.Cells.Item(99,6).Value = 0
$ws.Cells.Item(99,7).Value = 2
$ws.Cells.Item(99,8).Value = 125

# Row 104 (country -> Zambia)
$ws.Cells.Item(104,1).Value = "Zambia"
$ws.Cells.Item(104,2).Value = 3583
$ws.Cells.Item(104,3).Value = 197
$ws.Cells.Item(104,4).Value = 1677
$ws.Cells.Item(104,5).Value = 1778
$ws.Cells.Item(104,6).Value = 0
$ws.Cells.Item(104,7).Value = 0
$ws.Cells.Item(104,8).Value = 128

# Row 105 (country -> Nicaragua)
$ws.Cells.Item(105,1).Value = "Nicaragua"
$ws.Cells.Item(105,2).Value = 3439
$ws.Cells.Item(105,3).Value = 0
$ws.Cells.Item(105,4).Value = 2492
$ws.Cells.Item(105,5).Value = 839
$ws.Cells.Item(105,6).Value = 0
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = 108

# Row 128
$ws.Cells.Item(128,2).Value = 1840
$ws.Cells.Item(128,3).Value = 1
$ws.Cells.Item(128,4).Value = 1822
$ws.Cells.Item(128,5).Value = 8
$ws.Cells.Item(128,6).Value = 0
$ws.Cells.Item(128,7).Value = 0
$ws.Cells.Item(128,8).Value = 10

# Row 156
$ws.Cells.Item(156,2).Value = 679
$ws.Cells.Item(156,3).Value = 2
$ws.Cells.Item(156,4).Value = 665
$ws.Cells.Item(156,5).Value = 5
$ws.Cells.Item(156,6).Value = 0
$ws.Cells.Item(156,7).Value = 0
$ws.Cells.Item(156,8).Value = 9

# Row 162
$ws.Cells.Item(162,2).Value = 408
$ws.Cells.Item(162,3).Value = 7
$ws.Cells.Item(162,4).Value = 365
$ws.Cells.Item(162,5).Value = 43
$ws.Cells.Item(162,6).Value = 0
$ws.Cells.Item(162,7).Value = 0
$ws.Cells.Item(162,8).Value = 0

# Row 210 (country -> Islas Malvinas)
$ws.Cells.Item(210,1).Value = "Islas Malvinas"

# Row 211 (country -> Groenlandia)
$ws.Cells.Item(211,1).Value = "Groenlandia"
